# Workbook: Hortaliza, Mercado Mayorista Lo Valledor de Santiago - Acelga
# Weekly update: insert 3 new data rows (Extra/Primera/Segunda) for the new
# reporting date 44474 (2021-10-05) right after the header/constant block,
# pushing the existing data rows 436:500 down to 439:503.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at row 436; existing rows 436:500 shift down to 439:503
$ws.Rows("436:438").Insert()

# Common (constant) values shared by every data row in this sheet
$mercadoId   = 6
$mercado     = "Mercado Mayorista Lo Valledor de Santiago"
$region      = "Metropolitana"
$codreg      = 13
$categoriaId = 100112009
$categoria   = "Acelga"
$variedad    = "Sin especificar"
$unidad      = "`$/docena de atados"
$origen      = "Región Metropolitana"
$kgUnidades  = 3
$clasif      = "Hortaliza"
$fecha       = 44474

# New row data: Calidad, Volumen, PrecioMin, PrecioMax, PrecioProm, PrecioKg
$newRows = @(
    @("Extra",   150, 12000, 12000, 12000, 4000),
    @("Primera", 190, 10000, 10000, 10000, 3333),
    @("Segunda", 130,  8000,  8000,  8000, 2667)
)

$r = 436
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value  = $mercadoId
    $ws.Cells.Item($r, 2).Value  = $mercado
    $ws.Cells.Item($r, 3).Value  = $region
    $ws.Cells.Item($r, 4).Value  = $fecha
    $ws.Cells.Item($r, 5).Value  = $codreg
    $ws.Cells.Item($r, 6).Value  = $categoriaId
    $ws.Cells.Item($r, 7).Value  = $categoria
    $ws.Cells.Item($r, 8).Value  = $variedad
    $ws.Cells.Item($r, 9).Value  = $row[0]
    $ws.Cells.Item($r, 10).Value = $row[1]
    $ws.Cells.Item($r, 11).Value = $row[2]
    $ws.Cells.Item($r, 12).Value = $row[3]
    $ws.Cells.Item($r, 13).Value = $row[4]
    $ws.Cells.Item($r, 14).Value = $unidad
    $ws.Cells.Item($r, 15).Value = $origen
    $ws.Cells.Item($r, 16).Value = $row[5]
    $ws.Cells.Item($r, 17).Value = $kgUnidades
    $ws.Cells.Item($r, 18).Value = $clasif
    $r = $r + 1
}

Write-Output "done"
